# Fruta / hortaliza, semanal
#
# A new weekly record is inserted into the "Papa" (potato) price table for
# Feria Lagunitas de Puerto Montt, right after the existing row for date
# 44299 (row 316). All the rows that used to be 316..335 shift down by one
# (to 317..336); their contents are untouched by this edit.
#
# New row 316 values:
#   Fecha (D)                         44610
#   Variedad (H)                      Patagonia
#   Calidad (I)                       1a nueva(o)
#   Volumen (J)                       600
#   Precio minimo (K)                 6000
#   Precio maximo (L)                 7000
#   Precio promedio ponderado (M)     6500
#   Precio $/Kg (P)                   260
# (A, B, C, E, F, G, N, O, Q, R keep the same constant values used by every
#  other row of this sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing row from 316 downward one row, leaving a blank row 316
# (this is exactly what Excel's own "Insert Row" does, and it carries the
# row's formatting - e.g. the date style on column D - down with it).
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row.
$ws.Cells.Item(316, 1).Value  = 4
$ws.Cells.Item(316, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(316, 3).Value  = "Los Lagos"
$ws.Cells.Item(316, 4).Value  = 44610
$ws.Cells.Item(316, 5).Value  = 10
$ws.Cells.Item(316, 6).Value  = 100114001
$ws.Cells.Item(316, 7).Value  = "Papa"
$ws.Cells.Item(316, 8).Value  = "Patagonia"
$ws.Cells.Item(316, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(316, 10).Value = 600
$ws.Cells.Item(316, 11).Value = 6000
$ws.Cells.Item(316, 12).Value = 7000
$ws.Cells.Item(316, 13).Value = 6500
$ws.Cells.Item(316, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(316, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(316, 16).Value = 260
$ws.Cells.Item(316, 17).Value = 25
$ws.Cells.Item(316, 18).Value = "Hortaliza"
